# Apply the "missing_data" re-sampling edit described by the diff:
#  - Two whole data rows ("RM 232" and the original "SC 92" row) are
#    removed entirely, shifting every following row up.
#  - After the shift, a handful of individual cells in columns C/D/F
#    flip between "missing" (blank) and a concrete numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that disappear from the table -------------------
# Row 26 ("RM 232") is deleted first; what was row 28 ("SC 92") becomes
# row 27 once row 26 is gone, so we delete row 27 next.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Clear cells that become missing in the final layout ------------------
$cellsToClear = @("D3", "D4", "F5", "F12", "D13", "F18", "F19", "D25", "C29", "F29", "F33")
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# --- Fill cells that become populated (previously missing) in the final layout --
$ws.Range("D2").Value = -13.5
$ws.Range("F8").Value = 17.05
$ws.Range("F10").Value = 16.43
$ws.Range("D11").Value = -15.5
$ws.Range("F15").Value = 16.2
$ws.Range("D21").Value = -14.3
$ws.Range("F25").Value = 16.6
$ws.Range("F27").Value = 17
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
